$wb = $excel.ActiveWorkbook

# --- Sources sheet: update rate values ---
$sources = $wb.Worksheets.Item("Sources")
$sources.Range("C2").Value = 1000
$sources.Range("C3").Value = 2000

# --- Transformers sheet: swap loss-factor values between rows 5 and 6 ---
$transformers = $wb.Worksheets.Item("Transformers")
$transformers.Range("G5").Value = 0.2
$transformers.Range("G6").Value = 0.1

# --- CO2Locations sheet: set every capacity value (S2:S167) to 100 ---
$co2 = $wb.Worksheets.Item("CO2Locations")
$co2.Range("S2:S167").Value = 100

# --- Update the selection/active-cell bookmarks left in each sheet ---
# Sinks: selection moves from E2 to D3
$sinks = $wb.Worksheets.Item("Sinks")
[void]$sinks.Range("D3").Select()

# CO2Locations: selection moves from S6 to S2
[void]$co2.Range("S2").Select()

# Transformers: selection moves from G6 to G5, and it is no longer the tab shown on open
[void]$transformers.Range("G5").Select()

# Sources: selection moves from C5 to C3, and it becomes the tab shown on open (last
# Select() wins for which sheet is "active" / tabSelected on reopen)
[void]$sources.Range("C3").Select()
